$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage copies of source rows into a scratch area (rows 500+) to avoid overlap issues
# while permuting row contents within each block.

$ws.Range("A38:AY38").Copy() | Out-Null
$ws.Range("A500:AY500").PasteSpecial(-4163) | Out-Null
$ws.Range("A39:AY39").Copy() | Out-Null
$ws.Range("A501:AY501").PasteSpecial(-4163) | Out-Null
$ws.Range("A35:AY35").Copy() | Out-Null
$ws.Range("A502:AY502").PasteSpecial(-4163) | Out-Null
$ws.Range("A36:AY36").Copy() | Out-Null
$ws.Range("A503:AY503").PasteSpecial(-4163) | Out-Null
$ws.Range("A37:AY37").Copy() | Out-Null
$ws.Range("A504:AY504").PasteSpecial(-4163) | Out-Null
$ws.Range("A54:AY54").Copy() | Out-Null
$ws.Range("A505:AY505").PasteSpecial(-4163) | Out-Null
$ws.Range("A55:AY55").Copy() | Out-Null
$ws.Range("A506:AY506").PasteSpecial(-4163) | Out-Null
$ws.Range("A56:AY56").Copy() | Out-Null
$ws.Range("A507:AY507").PasteSpecial(-4163) | Out-Null
$ws.Range("A53:AY53").Copy() | Out-Null
$ws.Range("A508:AY508").PasteSpecial(-4163) | Out-Null
$ws.Range("A66:AY66").Copy() | Out-Null
$ws.Range("A509:AY509").PasteSpecial(-4163) | Out-Null
$ws.Range("A63:AY63").Copy() | Out-Null
$ws.Range("A510:AY510").PasteSpecial(-4163) | Out-Null
$ws.Range("A64:AY64").Copy() | Out-Null
$ws.Range("A511:AY511").PasteSpecial(-4163) | Out-Null
$ws.Range("A65:AY65").Copy() | Out-Null
$ws.Range("A512:AY512").PasteSpecial(-4163) | Out-Null

# Now copy staged rows into their final target row positions
$ws.Range("A500:AY500").Copy() | Out-Null
$ws.Range("A35:AY35").PasteSpecial(-4163) | Out-Null
$ws.Range("A501:AY501").Copy() | Out-Null
$ws.Range("A36:AY36").PasteSpecial(-4163) | Out-Null
$ws.Range("A502:AY502").Copy() | Out-Null
$ws.Range("A37:AY37").PasteSpecial(-4163) | Out-Null
$ws.Range("A503:AY503").Copy() | Out-Null
$ws.Range("A38:AY38").PasteSpecial(-4163) | Out-Null
$ws.Range("A504:AY504").Copy() | Out-Null
$ws.Range("A39:AY39").PasteSpecial(-4163) | Out-Null
$ws.Range("A505:AY505").Copy() | Out-Null
$ws.Range("A53:AY53").PasteSpecial(-4163) | Out-Null
$ws.Range("A506:AY506").Copy() | Out-Null
$ws.Range("A54:AY54").PasteSpecial(-4163) | Out-Null
$ws.Range("A507:AY507").Copy() | Out-Null
$ws.Range("A55:AY55").PasteSpecial(-4163) | Out-Null
$ws.Range("A508:AY508").Copy() | Out-Null
$ws.Range("A56:AY56").PasteSpecial(-4163) | Out-Null
$ws.Range("A509:AY509").Copy() | Out-Null
$ws.Range("A63:AY63").PasteSpecial(-4163) | Out-Null
$ws.Range("A510:AY510").Copy() | Out-Null
$ws.Range("A64:AY64").PasteSpecial(-4163) | Out-Null
$ws.Range("A511:AY511").Copy() | Out-Null
$ws.Range("A65:AY65").PasteSpecial(-4163) | Out-Null
$ws.Range("A512:AY512").Copy() | Out-Null
$ws.Range("A66:AY66").PasteSpecial(-4163) | Out-Null

# Clean up staging area
$ws.Range("A500:AY500").ClearContents() | Out-Null
$ws.Range("A501:AY501").ClearContents() | Out-Null
$ws.Range("A502:AY502").ClearContents() | Out-Null
$ws.Range("A503:AY503").ClearContents() | Out-Null
$ws.Range("A504:AY504").ClearContents() | Out-Null
$ws.Range("A505:AY505").ClearContents() | Out-Null
$ws.Range("A506:AY506").ClearContents() | Out-Null
$ws.Range("A507:AY507").ClearContents() | Out-Null
$ws.Range("A508:AY508").ClearContents() | Out-Null
$ws.Range("A509:AY509").ClearContents() | Out-Null
$ws.Range("A510:AY510").ClearContents() | Out-Null
$ws.Range("A511:AY511").ClearContents() | Out-Null
$ws.Range("A512:AY512").ClearContents() | Out-Null

$excel.CutCopyMode = 0
